$wb = $excel.ActiveWorkbook

# =====================================================================
# 1) Insert a brand-new "2022-Q4" worksheet right after "总计", i.e.
#    right before the existing "2021-Q4" sheet (which - together with
#    "2021-Q3" - simply shifts one tab to the right).
# =====================================================================
$existingQ4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($existingQ4)
$newSheet.Name = "2022-Q4"

# Reuse the sibling quarter sheet's look (header style in B1:H1 plus the
# bold/centered/bordered "index" style used down column A).
$srcFmt = $wb.Worksheets.Item("2021-Q4").Range("B1:H3")
$srcFmt.Copy()
$newSheet.Range("B1").PasteSpecial(-4122)

$srcA = $wb.Worksheets.Item("2021-Q4").Range("A2:A3")
$srcA.Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A4:A5").PasteSpecial(-4122)

# Helper: write a value as plain TEXT (so numeric-looking strings such as
# "012866" or "3.24" keep their leading zeros / do not get coerced into
# numbers), without leaving a leftover "entered with a leading quote"
# marker behind - re-pasting an untouched blank cell's (default) format
# over the cell clears that marker while the text value itself survives.
$blankFmt = $newSheet.Range("Z1")
function Set-PlainText($range, [string]$text) {
    $range.Value = "'" + $text
    $blankFmt.Copy()
    $range.PasteSpecial(-4122)
}

# ---- Header row ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- Data rows (fund codes / sizes / ratios are stored as text, like the
#      other quarter tabs; H is a genuine number) ----
$newSheet.Range("A2").Value = 0
Set-PlainText $newSheet.Range("B2") "012866"
$newSheet.Range("C2").Value = "易方达标普生物科技指数（QDII-LOF）人民币 C"
Set-PlainText $newSheet.Range("D2") "3.24"
Set-PlainText $newSheet.Range("E2") "94.17"
Set-PlainText $newSheet.Range("F2") "1.16"
Set-PlainText $newSheet.Range("G2") "0.0376"
$newSheet.Range("H2").Value = 3

$newSheet.Range("A3").Value = 1
Set-PlainText $newSheet.Range("B3") "161127"
$newSheet.Range("C3").Value = "易方达标普生物科技指数（QDII-LOF）人民币"
Set-PlainText $newSheet.Range("D3") "3.24"
Set-PlainText $newSheet.Range("E3") "94.17"
Set-PlainText $newSheet.Range("F3") "1.16"
Set-PlainText $newSheet.Range("G3") "0.0376"
$newSheet.Range("H3").Value = 3

$newSheet.Range("A4").Value = 2
Set-PlainText $newSheet.Range("B4") "003720"
$newSheet.Range("C4").Value = "易方达标普生物科技指数（QDII-LOF）美元A"
Set-PlainText $newSheet.Range("D4") "3.12"
Set-PlainText $newSheet.Range("E4") "94.17"
Set-PlainText $newSheet.Range("F4") "1.16"
Set-PlainText $newSheet.Range("G4") "0.0362"
$newSheet.Range("H4").Value = 3

$newSheet.Range("A5").Value = 3
Set-PlainText $newSheet.Range("B5") "012867"
$newSheet.Range("C5").Value = "易方达标普生物科技指数（QDII-LOF）美元 C"
Set-PlainText $newSheet.Range("D5") "0.12"
Set-PlainText $newSheet.Range("E5") "94.17"
Set-PlainText $newSheet.Range("F5") "1.16"
Set-PlainText $newSheet.Range("G5") "0.0014"
$newSheet.Range("H5").Value = 3

# =====================================================================
# 2) Update the "总计" summary sheet: the new 2022-Q4 quarter becomes the
#    first data row, pushing 2021-Q4 / 2021-Q3 down by one row each.
# =====================================================================
$summary = $wb.Worksheets.Item("总计")

# Row 4 is brand new - give A4 the same "index" style used by A2 / A3
# before filling in the (old row-3 / 2021-Q3) values.
$summary.Range("A4").Value = 2
$summary.Range("A3").Copy()
$summary.Range("A4").PasteSpecial(-4122)
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q3"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.03

# Row 3 now holds what used to be row 2 (2021-Q4)
$summary.Range("B3").Value = "2021-Q4"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.05

# Row 2 now holds the brand new 2022-Q4 figures
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.11

# =====================================================================
# 3) Restore the original active tab: "2021-Q3" was the selected sheet
#    before the edit and stays selected afterwards.
# =====================================================================
$wb.Worksheets.Item("2021-Q3").Activate()
